$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Antal" column (I) holds text-looking numbers ("10"/"20") that must
# stay text, not get auto-converted to numeric by Excel's input parser.
# Pre-formatting the cell as Text before assigning keeps it a string.
$ws.Range("I7:I9").NumberFormat = "@"

# Row 7: was Id 112092161 / Antal "10" -> becomes Id 112092130 / Antal blank
$ws.Range("A7").Value = 112092130
$ws.Range("I7").Value = ""
$ws.Range("Q7").Value = 584352
$ws.Range("R7").Value = 7048232

# Row 8: was Id 112092586 / Antal "20" -> becomes Id 112092161 / Antal "10"
$ws.Range("A8").Value = 112092161
$ws.Range("I8").Value = "10"
$ws.Range("Q8").Value = 584330
$ws.Range("R8").Value = 7048274
$ws.Range("Z8").Value = "17:22"
$ws.Range("AB8").Value = "17:22"

# Row 9: was Id 112092130 / Antal blank -> becomes Id 112092586 / Antal "20"
$ws.Range("A9").Value = 112092586
$ws.Range("I9").Value = "20"
$ws.Range("Q9").Value = 584401
$ws.Range("R9").Value = 7048357
$ws.Range("Z9").Value = "17:46"
$ws.Range("AB9").Value = "17:46"

# Row 10: Id/Antal unchanged, only Ost/Nord coordinates get rounded
$ws.Range("Q10").Value = 584346
$ws.Range("R10").Value = 7048207
